$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 27
$ws.Range("H27").Value = 4000
$ws.Range("J27").Value = 4000
$ws.Range("L27").Value = 12000
$ws.Range("N27").Value = -12202
# Row 41
$ws.Range("H41").Value = 289
$ws.Range("I41").Value = 343.75
$ws.Range("J41").Value = 179.5
$ws.Range("K41").Value = 343.75
$ws.Range("L41").Value = 179.5
$ws.Range("M41").Value = 96.25
$ws.Range("N41").Value = -1059.5
# Row 81
$ws.Range("H81").Value = 90000
$ws.Range("J81").Value = 90000
$ws.Range("L81").Value = 90000
$ws.Range("N81").Value = -91996
# Row 84
$ws.Range("H84").Value = 90000
$ws.Range("J84").Value = 90000
$ws.Range("L84").Value = 270000
$ws.Range("N84").Value = -279984
# Row 106
$ws.Range("H106").Value = 3031
$ws.Range("I106").Value = 2468.4285
$ws.Range("K106").Value = 2468.4285
$ws.Range("M106").Value = -1837.4285
# Row 121
$ws.Range("H121").Value = 8583.788
$ws.Range("J121").Value = 8768.23
$ws.Range("L121").Value = 26304.69
$ws.Range("N121").Value = -29798.69
# Row 131
$ws.Range("H131").Value = 169469.83
$ws.Range("I131").Value = 253023.5
$ws.Range("J131").Value = 2362.5
$ws.Range("K131").Value = 759070.5
$ws.Range("L131").Value = 7087.5
$ws.Range("M131").Value = -754030.5
$ws.Range("N131").Value = -17167.5
# Row 132
$ws.Range("H132").Value = 8527.267
$ws.Range("I132").Value = 10701
$ws.Range("K132").Value = 32103
$ws.Range("M132").Value = -29573
# Row 137
$ws.Range("H137").Value = 2246.125
$ws.Range("I137").Value = 1793.8
$ws.Range("K137").Value = 5381.4
$ws.Range("M137").Value = -2831.4
# Row 138
$ws.Range("H138").Value = 2817.5652
$ws.Range("J138").Value = 3072.5806
$ws.Range("L138").Value = 9217.7418
$ws.Range("N138").Value = -19497.7418
# Row 141
$ws.Range("H141").Value = 5439.1904
$ws.Range("I141").Value = 4681.8887
$ws.Range("K141").Value = 14045.6661
$ws.Range("M141").Value = -8865.666100000002

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 5383443
$ws.Range("I61").Value = 8777408
$ws.Range("J61").Value = 9666
$ws.Range("K61").Value = 8777408
$ws.Range("L61").Value = 9666
$ws.Range("M61").Value = -8777196
$ws.Range("N61").Value = -10090
# Row 74
$ws.Range("H74").Value = 3716.6538
$ws.Range("I74").Value = 2776.375
$ws.Range("J74").Value = 15000
$ws.Range("K74").Value = 2776.375
$ws.Range("L74").Value = 15000
$ws.Range("M74").Value = -1902.375
$ws.Range("N74").Value = -16748
# Row 76
$ws.Range("H76").Value = 17952.666
$ws.Range("J76").Value = 17952.666
$ws.Range("L76").Value = 17952.666
$ws.Range("N76").Value = -18628.666
# Row 77
$ws.Range("H77").Value = 3716.6538
$ws.Range("I77").Value = 2776.375
$ws.Range("J77").Value = 15000
$ws.Range("K77").Value = 13881.875
$ws.Range("L77").Value = 75000
$ws.Range("M77").Value = -9513.875
$ws.Range("N77").Value = -83736
# Row 79
$ws.Range("H79").Value = 17952.666
$ws.Range("J79").Value = 17952.666
$ws.Range("L79").Value = 17952.666
$ws.Range("N79").Value = -20292.666
# Row 132
$ws.Range("H132").Value = 6120.724
$ws.Range("I132").Value = 5100.4
$ws.Range("K132").Value = 15301.2
$ws.Range("M132").Value = -12771.2
# Row 136
$ws.Range("H136").Value = 5383443
$ws.Range("I136").Value = 8777408
$ws.Range("J136").Value = 9666
$ws.Range("K136").Value = 26332224
$ws.Range("L136").Value = 28998
$ws.Range("M136").Value = -26329674
$ws.Range("N136").Value = -34098

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 46
$ws.Range("H46").Value = 29999.8
$ws.Range("J46").Value = 29999.8
$ws.Range("L46").Value = 29999.8
$ws.Range("N46").Value = -30595.8
# Row 134
$ws.Range("H134").Value = 3796.2974
$ws.Range("I134").Value = 3846.1943
$ws.Range("K134").Value = 11538.5829
$ws.Range("M134").Value = -9003.582900000001

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4324.684
$ws.Range("I31").Value = 3063.9656
$ws.Range("K31").Value = 3063.9656
$ws.Range("M31").Value = -2768.9656
# Row 32
$ws.Range("H32").Value = 3200
$ws.Range("I32").Value = 3200
$ws.Range("K32").Value = 3200
$ws.Range("M32").Value = -2884
# Row 34
$ws.Range("H34").Value = 4324.684
$ws.Range("I34").Value = 3063.9656
$ws.Range("K34").Value = 3063.9656
$ws.Range("M34").Value = -2861.9656
# Row 39
$ws.Range("H39").Value = 449.5
$ws.Range("I39").Value = 449.5
$ws.Range("K39").Value = 449.5
$ws.Range("M39").Value = -58.5
# Row 49
$ws.Range("H49").Value = 449.5
$ws.Range("I49").Value = 449.5
$ws.Range("K49").Value = 449.5
$ws.Range("M49").Value = -267.5
# Row 94
$ws.Range("H94").Value = 4829.2856
$ws.Range("J94").Value = 4829.2856
$ws.Range("L94").Value = 4829.2856
$ws.Range("N94").Value = -5731.2856
# Row 132
$ws.Range("H132").Value = 1944.8125
$ws.Range("I132").Value = 1593.9166
$ws.Range("J132").Value = 2997.5
$ws.Range("K132").Value = 4781.7498
$ws.Range("L132").Value = 8992.5
$ws.Range("M132").Value = -2251.7498
$ws.Range("N132").Value = -14052.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 4403.875
$ws.Range("I80").Value = 3922
$ws.Range("K80").Value = 3922
$ws.Range("M80").Value = -2924
# Row 83
$ws.Range("H83").Value = 4403.875
$ws.Range("I83").Value = 3922
$ws.Range("K83").Value = 19610
$ws.Range("M83").Value = -14618
# Row 107
$ws.Range("H107").Value = 1303960
$ws.Range("I107").Value = 2607449.2
$ws.Range("J107").Value = 470.7143
$ws.Range("K107").Value = 2607449.2
$ws.Range("L107").Value = 470.7143
$ws.Range("M107").Value = -2605529.2
$ws.Range("N107").Value = -4310.7143
# Row 122
$ws.Range("H122").Value = 2728.2856
$ws.Range("I122").Value = 2728.2856
$ws.Range("K122").Value = 8184.8568
$ws.Range("M122").Value = -5734.8568

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 3472.5
# Row 24
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
# Row 27
$ws.Range("H27").Value = 3472.5
# Row 68
$ws.Range("H68").Value = 7958.2354
$ws.Range("J68").Value = 6941.3335
$ws.Range("L68").Value = 6941.3335
$ws.Range("N68").Value = -8439.333500000001
# Row 71
$ws.Range("H71").Value = 7958.2354
$ws.Range("J71").Value = 6941.3335
$ws.Range("L71").Value = 34706.6675
$ws.Range("N71").Value = -42194.6675
# Row 93
$ws.Range("H93").Value = 17702.846
$ws.Range("I93").Value = 1014
$ws.Range("K93").Value = 1014
$ws.Range("M93").Value = 234
# Row 107
$ws.Range("H107").Value = 17029.715
$ws.Range("I107").Value = 17029.715
$ws.Range("K107").Value = 17029.715
$ws.Range("M107").Value = -15109.715

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 7
$ws.Range("H7").Value = 1449
$ws.Range("I7").Value = 1449
$ws.Range("K7").Value = 1449
$ws.Range("M7").Value = -1336
# Row 62
$ws.Range("H62").Value = 9998.799999999999
$ws.Range("J62").Value = 9998.799999999999
$ws.Range("L62").Value = 9998.799999999999
$ws.Range("N62").Value = -11246.8
# Row 65
$ws.Range("H65").Value = 9998.799999999999
$ws.Range("J65").Value = 9998.799999999999
$ws.Range("L65").Value = 49994
$ws.Range("N65").Value = -56234
# Row 80
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
# Row 83
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
# Row 122
$ws.Range("H122").Value = 3185.9524
$ws.Range("I122").Value = 2470.4167
$ws.Range("K122").Value = 7411.250100000001
$ws.Range("M122").Value = -4961.250100000001
# Row 132
$ws.Range("H132").Value = 5278.881
$ws.Range("I132").Value = 4665.4517
$ws.Range("K132").Value = 13996.3551
$ws.Range("M132").Value = -11466.3551
# Row 136
$ws.Range("H136").Value = 4997.575
$ws.Range("I136").Value = 3983.4243
$ws.Range("J136").Value = 9778.571
$ws.Range("K136").Value = 11950.2729
$ws.Range("L136").Value = 29335.713
$ws.Range("M136").Value = -9400.2729
$ws.Range("N136").Value = -34435.713

$wb.Save()
Write-Host "Applied all changes"